$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 440 (shifts old rows 440..474 down to 441..475)
$ws.Rows.Item(440).Insert()

# Populate the newly inserted row 440 with the new weekly price entry.
# (Mirrors the surrounding rows' fixed columns for this market/category.)
$ws.Cells.Item(440, 1).Value = 11
$ws.Cells.Item(440, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(440, 3).Value = "Bíobío"
$ws.Cells.Item(440, 4).Value = 45212
$ws.Cells.Item(440, 5).Value = 8
$ws.Cells.Item(440, 6).Value = 100112009
$ws.Cells.Item(440, 7).Value = "Acelga"
$ws.Cells.Item(440, 8).Value = "Sin especificar"
$ws.Cells.Item(440, 9).Value = "Primera"
$ws.Cells.Item(440, 10).Value = 180
$ws.Cells.Item(440, 11).Value = 600
$ws.Cells.Item(440, 12).Value = 650
$ws.Cells.Item(440, 13).Value = 628
$ws.Cells.Item(440, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(440, 15).Value = "Región de Ñuble"
$ws.Cells.Item(440, 16).Value = 628
$ws.Cells.Item(440, 17).Value = 1
$ws.Cells.Item(440, 18).Value = "Hortaliza"
